# Publish documentation 0.1.1 / ror 0.1.1
# - Bump the "Version" metadata row from 0.1.0 to 0.1.1
# - Bump the "Date" metadata row to the new publication date
# - Add a new "Context" metadata row for element:ContactPoint

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: Version / 0.1.0 -> 0.1.1
$ws.Cells.Item(3, 2).Value = "0.1.1"

# Row 8: Date -> new publish date
$ws.Cells.Item(8, 2).Value = "2023-06-02T12:02:38+02:00"

# New row 21: Context / element:ContactPoint (mirrors the existing
# Context row 20, so copy its formatting first)
$ws.Cells.Item(21, 1).Value = "Context"
$ws.Cells.Item(21, 2).Value = "element:ContactPoint"

$ws.Range("A20:B20").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
